# Insert a new data record (row) into the daily price log at row 425.
# Excel shifts the existing rows 425:549 down to 426:550 and extends the
# used range to A1:R550, exactly like a manual "Insert Sheet Rows".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(425).Insert()

# Populate the newly inserted row with the new record's values.
$ws.Cells.Item(425, 1).Value  = 4
$ws.Cells.Item(425, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(425, 3).Value  = "Los Lagos"
$ws.Cells.Item(425, 4).Value  = 45093
$ws.Cells.Item(425, 5).Value  = 10
$ws.Cells.Item(425, 6).Value  = 100112008
$ws.Cells.Item(425, 7).Value  = "Coliflor"
$ws.Cells.Item(425, 8).Value  = "Sin especificar"
$ws.Cells.Item(425, 9).Value  = "Primera"
$ws.Cells.Item(425, 10).Value = 1400
$ws.Cells.Item(425, 11).Value = 1500
$ws.Cells.Item(425, 12).Value = 1500
$ws.Cells.Item(425, 13).Value = 1500
$ws.Cells.Item(425, 14).Value = "$/unidad"
$ws.Cells.Item(425, 15).Value = "Región Metropolitana"
$ws.Cells.Item(425, 16).Value = 1500
$ws.Cells.Item(425, 17).Value = 1
$ws.Cells.Item(425, 18).Value = "Hortaliza"
